# Add three new columns (F, G, H) holding MAD-based outlier flags for the
# KNN, SVM and RF imputation methods, matching the existing worksheet
# layout/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
# Copy the style of the existing header cell (E1) onto the new header
# cells so the new headers are bold/centered/bordered like the rest.
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# --- Data rows (rows 2-17) -------------------------------------------------
# All cells default to FALSE, except F12 and F16 which are TRUE.
$trueRows = @(12, 16)

for ($r = 2; $r -le 17; $r++) {
    if ($trueRows -contains $r) {
        $ws.Cells.Item($r, 6).Value = $true
    } else {
        $ws.Cells.Item($r, 6).Value = $false
    }
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}
